$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Swap the data (columns B:AC) between several row pairs.
#    Column A (the sequential row id) stays put; only the match
#    data that follows it moves between the two rows.
# ---------------------------------------------------------------
$pairs = @(
    @(9, 10),
    @(49, 50),
    @(76, 77),
    @(87, 88)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")
    $v1 = $range1.Value()
    $v2 = $range2.Value()
    $range1.Value = $v2
    $range2.Value = $v1
}

# ---------------------------------------------------------------
# 2) Append two new match rows at the bottom of the sheet.
# ---------------------------------------------------------------

# Row 154 mirrors the formatting of the previous last row (153):
#  - column A uses the bold/bordered/centered style
#  - column E uses the date/time number format
$ws.Cells.Item(153, 1).Copy()
$ws.Cells.Item(154, 1).PasteSpecial(-4122)
$ws.Cells.Item(153, 5).Copy()
$ws.Cells.Item(154, 5).PasteSpecial(-4122)
$ws.Cells.Item(153, 1).Copy()
$ws.Cells.Item(155, 1).PasteSpecial(-4122)
$ws.Cells.Item(153, 5).Copy()
$ws.Cells.Item(155, 5).PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$data154 = @(152, 7952745, "Bosnia Herzegovina Premier Liga", "Bosnia  Herzegovina Premier Liga", 45395.35416666666, "Zvijezda 09", "GOSK Gabela", 1, 2, "A", 2.1, 3, 3.25, 2, 3.1, 3.3, -0.25, 1.8, 2, 2.25, 1.8, 2, -1, -1, 2.3, -1, 1, 0.8, -1)

$data155 = @(153, 7952458, "Bosnia Herzegovina Premier Liga", "Bosnia  Herzegovina Premier Liga", 45395.45833333334, "Zrinjski Mostar", "Velez Mostar", 1, 0, "H", 1.615, 3.4, 5, 1.533, 3.4, 5.75, -1, 2.025, 1.775, 2, 1.775, 2.025, 0.5329999999999999, -1, -1, 0, -0, -1, 1.025)

$row154 = New-Object 'object[,]' 1,29
for ($i = 0; $i -lt 29; $i++) { $row154[0,$i] = $data154[$i] }

$row155 = New-Object 'object[,]' 1,29
for ($i = 0; $i -lt 29; $i++) { $row155[0,$i] = $data155[$i] }

$ws.Range("A154:AC154").Value = $row154
$ws.Range("A155:AC155").Value = $row155
